# Generate Report for Handoff
# Adds two newly-ready-for-handoff files (308dd797-... and f2f4a233-...)
# to the Overview / zh-cn / de-de localization-status tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3) -- columns: File Name | Path And Name | Extension
#   | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add()
$wsOverview.Range("A4").Value2 = "308dd797-ae76-4838-90c4-a761eab4680c.md"
$wsOverview.Range("B4").Value2 = "e2e\308dd797-ae76-4838-90c4-a761eab4680c.md"
$wsOverview.Range("B4").Style = "Hyperlink"
$wsOverview.Range("C4").Value2 = ".md"
$wsOverview.Range("D4").Value2 = "'"
$wsOverview.Range("E4").Value2 = "Ready for handoff"
$wsOverview.Range("F4").Value2 = "Ready for handoff"
$wsOverview.Range("G4").Value2 = "2016-08-16 22:41:21"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/308dd797ae76483890c4a761eab4680cbuild/e2e/308dd797-ae76-4838-90c4-a761eab4680c.md", "", "", "e2e\308dd797-ae76-4838-90c4-a761eab4680c.md")

$loOverview.ListRows.Add()
$wsOverview.Range("A5").Value2 = "f2f4a233-5d05-4655-a557-132d5b5a60a6.md"
$wsOverview.Range("B5").Value2 = "e2e\f2f4a233-5d05-4655-a557-132d5b5a60a6.md"
$wsOverview.Range("B5").Style = "Hyperlink"
$wsOverview.Range("C5").Value2 = ".md"
$wsOverview.Range("D5").Value2 = "'"
$wsOverview.Range("E5").Value2 = "Ready for handoff"
$wsOverview.Range("F5").Value2 = "Ready for handoff"
$wsOverview.Range("G5").Value2 = "2016-08-16 22:41:21"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2f4a2335d054655a557132d5b5a60a6build/e2e/f2f4a233-5d05-4655-a557-132d5b5a60a6.md", "", "", "e2e\f2f4a233-5d05-4655-a557-132d5b5a60a6.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1) -- columns: Source File Name | File Extension | Status
#   | Source Path | Priority | Content Duplicate | Latest Handoff File
#   | Latest Handoff Datetime | Latest Target File | Latest Handback File
#   | Latest Handback DateTime | Reference Tokens | To be localized
#   | Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add()
$wsZhCn.Range("A4").Value2 = "308dd797-ae76-4838-90c4-a761eab4680c.md"
$wsZhCn.Range("A4").Style = "Hyperlink"
$wsZhCn.Range("B4").Value2 = ".md"
$wsZhCn.Range("C4").Value2 = "Ready for handoff"
$wsZhCn.Range("D4").Value2 = "e2e"
$wsZhCn.Range("E4").Value2 = "ht"
$wsZhCn.Range("F4").Value2 = "'False"
$wsZhCn.Range("G4").Value2 = "308dd797-ae76-4838-90c4-a761eab4680c.9f0420ef8ff9a4870ead7a5fb6a3b46025ef0979.zh-cn.xlf"
$wsZhCn.Range("H4").Value2 = "2016-08-16 22:41:16"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value2 = "'"
$wsZhCn.Range("J4").Value2 = "'"
$wsZhCn.Range("K4").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value2 = "'"
$wsZhCn.Range("M4").Value2 = "'True"
$wsZhCn.Range("N4").Value2 = "'"
$wsZhCn.Range("O4").Value2 = "'False"
$wsZhCn.Range("P4").Value2 = "'"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/308dd797ae76483890c4a761eab4680cbuild/e2e/308dd797-ae76-4838-90c4-a761eab4680c.md", "", "", "308dd797-ae76-4838-90c4-a761eab4680c.md")

$loZhCn.ListRows.Add()
$wsZhCn.Range("A5").Value2 = "f2f4a233-5d05-4655-a557-132d5b5a60a6.md"
$wsZhCn.Range("A5").Style = "Hyperlink"
$wsZhCn.Range("B5").Value2 = ".md"
$wsZhCn.Range("C5").Value2 = "Ready for handoff"
$wsZhCn.Range("D5").Value2 = "e2e"
$wsZhCn.Range("E5").Value2 = "ht"
$wsZhCn.Range("F5").Value2 = "'False"
$wsZhCn.Range("G5").Value2 = "f2f4a233-5d05-4655-a557-132d5b5a60a6.a2e4330169c4ed4312e862f966b4a1aadd0dc6c4.zh-cn.xlf"
$wsZhCn.Range("H5").Value2 = "2016-08-16 22:41:16"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I5").Value2 = "'"
$wsZhCn.Range("J5").Value2 = "'"
$wsZhCn.Range("K5").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L5").Value2 = "'"
$wsZhCn.Range("M5").Value2 = "'True"
$wsZhCn.Range("N5").Value2 = "'"
$wsZhCn.Range("O5").Value2 = "'False"
$wsZhCn.Range("P5").Value2 = "'"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2f4a2335d054655a557132d5b5a60a6build/e2e/f2f4a233-5d05-4655-a557-132d5b5a60a6.md", "", "", "f2f4a233-5d05-4655-a557-132d5b5a60a6.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2) -- same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add()
$wsDeDe.Range("A4").Value2 = "308dd797-ae76-4838-90c4-a761eab4680c.md"
$wsDeDe.Range("A4").Style = "Hyperlink"
$wsDeDe.Range("B4").Value2 = ".md"
$wsDeDe.Range("C4").Value2 = "Ready for handoff"
$wsDeDe.Range("D4").Value2 = "e2e"
$wsDeDe.Range("E4").Value2 = "ht"
$wsDeDe.Range("F4").Value2 = "'False"
$wsDeDe.Range("G4").Value2 = "308dd797-ae76-4838-90c4-a761eab4680c.9f0420ef8ff9a4870ead7a5fb6a3b46025ef0979.de-de.xlf"
$wsDeDe.Range("H4").Value2 = "2016-08-16 22:41:21"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value2 = "'"
$wsDeDe.Range("J4").Value2 = "'"
$wsDeDe.Range("K4").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value2 = "'"
$wsDeDe.Range("M4").Value2 = "'True"
$wsDeDe.Range("N4").Value2 = "'"
$wsDeDe.Range("O4").Value2 = "'False"
$wsDeDe.Range("P4").Value2 = "'"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/308dd797ae76483890c4a761eab4680cbuild/e2e/308dd797-ae76-4838-90c4-a761eab4680c.md", "", "", "308dd797-ae76-4838-90c4-a761eab4680c.md")

$loDeDe.ListRows.Add()
$wsDeDe.Range("A5").Value2 = "f2f4a233-5d05-4655-a557-132d5b5a60a6.md"
$wsDeDe.Range("A5").Style = "Hyperlink"
$wsDeDe.Range("B5").Value2 = ".md"
$wsDeDe.Range("C5").Value2 = "Ready for handoff"
$wsDeDe.Range("D5").Value2 = "e2e"
$wsDeDe.Range("E5").Value2 = "ht"
$wsDeDe.Range("F5").Value2 = "'False"
$wsDeDe.Range("G5").Value2 = "f2f4a233-5d05-4655-a557-132d5b5a60a6.a2e4330169c4ed4312e862f966b4a1aadd0dc6c4.de-de.xlf"
$wsDeDe.Range("H5").Value2 = "2016-08-16 22:41:21"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I5").Value2 = "'"
$wsDeDe.Range("J5").Value2 = "'"
$wsDeDe.Range("K5").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L5").Value2 = "'"
$wsDeDe.Range("M5").Value2 = "'True"
$wsDeDe.Range("N5").Value2 = "'"
$wsDeDe.Range("O5").Value2 = "'False"
$wsDeDe.Range("P5").Value2 = "'"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2f4a2335d054655a557132d5b5a60a6build/e2e/f2f4a233-5d05-4655-a557-132d5b5a60a6.md", "", "", "f2f4a233-5d05-4655-a557-132d5b5a60a6.md")
